# Weekly update of the "Cebollín" price sheet.
# Rewrites the data rows (prices shuffled/refreshed for the new week) and
# appends one new row (row 15) with the same shape as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the date-column number format from an existing data row so the new
# row 15 date cell (and any row whose style may have been lost) matches.
$dateFmt = $ws.Cells.Item(2, 4).NumberFormat

function Set-Row {
    param($Row, $Dt, $Vol, $Pmin, $Pmax, $Pprom, $Unidad, $Origen, $Pkg, $KgU)

    $ws.Cells.Item($Row, 4).Value = $Dt
    $ws.Cells.Item($Row, 4).NumberFormat = $dateFmt
    $ws.Cells.Item($Row, 10).Value = $Vol
    $ws.Cells.Item($Row, 11).Value = $Pmin
    $ws.Cells.Item($Row, 12).Value = $Pmax
    $ws.Cells.Item($Row, 13).Value = $Pprom
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $Pkg
    $ws.Cells.Item($Row, 17).Value = $KgU
}

# Row 2
Set-Row 2 44704 100 6000 6500 6250 "$/paquete 36 unidades" "Región Metropolitana" 174 36

# Row 3 unchanged

# Row 4
Set-Row 4 44760 120 8000 8000 8000 "$/docena de atados" "Región Metropolitana" 2667 3

# Row 5
Set-Row 5 44166 70 3500 4000 3679 "$/paquete 36 unidades" "Región Metropolitana" 102 36

# Row 6
Set-Row 6 44160 43 3500 4000 3709 "$/paquete 36 unidades" "Región Metropolitana" 103 36

# Row 7 unchanged

# Row 8
Set-Row 8 44209 150 3500 4000 3767 "$/paquete 2 kilos" "Provincia de Diguillín" 1884 2

# Row 9
Set-Row 9 44225 80 3400 3700 3550 "$/paquete 2 kilos" "Provincia de Diguillín" 1775 2

# Row 10
Set-Row 10 44664 200 8000 8500 8250 "$/paquete 36 unidades" "Región Metropolitana" 229 36

# Row 11
Set-Row 11 44215 140 3500 4000 3768 "$/paquete 2 kilos" "Provincia de Diguillín" 1884 2

# Row 12
Set-Row 12 44223 80 3500 3800 3688 "$/paquete 2 kilos" "Provincia de Diguillín" 1844 2

# Row 13
Set-Row 13 44210 105 3500 4000 3714 "$/paquete 2 kilos" "Provincia de Diguillín" 1857 2

# Row 14
Set-Row 14 44161 50 2800 3000 2900 "$/paquete 2 kilos" "Provincia de Diguillín" 1450 2

# Row 15 (new row; constant columns copied from row 14, which has the same
# Mercado/Region/Categoria/etc. as every other row in this sheet)
$ws.Cells.Item(15, 1).Value = $ws.Cells.Item(14, 1).Value2
$ws.Cells.Item(15, 2).Value = $ws.Cells.Item(14, 2).Value2
$ws.Cells.Item(15, 3).Value = $ws.Cells.Item(14, 3).Value2
$ws.Cells.Item(15, 5).Value = $ws.Cells.Item(14, 5).Value2
$ws.Cells.Item(15, 6).Value = $ws.Cells.Item(14, 6).Value2
$ws.Cells.Item(15, 7).Value = $ws.Cells.Item(14, 7).Value2
$ws.Cells.Item(15, 8).Value = $ws.Cells.Item(14, 8).Value2
$ws.Cells.Item(15, 9).Value = $ws.Cells.Item(14, 9).Value2
$ws.Cells.Item(15, 18).Value = $ws.Cells.Item(14, 18).Value2

Set-Row 15 44662 200 8000 8500 8250 "$/paquete 36 unidades" "Región Metropolitana" 229 36
